$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q4", mirroring the layout
#    of the "2021-Q4" sheet (same headers/style) with the new quarter's data.
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $wsQ4)
$newSheet.Name = "2022-Q1"

# Copy formatting + values for the header row and the B2:H2 data cells from
# "2021-Q4" (identical headers), then overwrite with the 2022-Q1 values.
$wsQ4.Range("B1:H2").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$wsQ4.Range("B1:H2").Copy()
$newSheet.Range("B1").PasteSpecial(-4163)   # xlPasteValues (keeps text typing)

$wsQ4.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$newSheet.Range("A2").PasteSpecial(-4163)   # xlPasteValues

$newSheet.Range("B2").Value = "'007280"
$newSheet.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$newSheet.Range("D2").Value = "'1.35"
$newSheet.Range("E2").Value = "'88.71"
$newSheet.Range("F2").Value = "'1.98"
$newSheet.Range("G2").Value = "'0.0267"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert the new 2022-Q1 row at the top
#    of the data, pushing the older quarters down one row, and append the
#    2020-Q4 row that now falls off the bottom.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Give the new A7 index cell the same style as the existing index column.
$wsTotal.Range("A6").Copy()
$wsTotal.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.03

$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.04

$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.04

$wsTotal.Range("B5").Value = "2021-Q2"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 0.05

$wsTotal.Range("B6").Value = "2021-Q1"
$wsTotal.Range("C6").Value = 1
$wsTotal.Range("D6").Value = 0.04

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2020-Q4"
$wsTotal.Range("C7").Value = 1
$wsTotal.Range("D7").Value = 0.05
